$d = $word.ActiveDocument

# Disable smart-quote autocorrect so literal straight quotes/apostrophes
# inserted via Find & Replace are not silently "curlified".
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# ---------------------------------------------------------------------------
# 1) Title heading text: "Nov 2017" -> "August 2018"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ashley Hindmarsh - Curriculum Vitae - Nov 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ashley Hindmarsh - Curriculum Vitae - August 2018", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Rename the matching bookmark (keep it a zero-length bookmark located
#    right before the heading run, just like the original).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("ashley-hindmarsh---curriculum-vitae---nov-2017")
$bm.Delete()

$tmp = $d.Range(0, 0)
$tmp.InsertBefore("X")
$newBmRange = $d.Range(0, 1)
$newBmRange.Bookmarks.Add("ashley-hindmarsh---curriculum-vitae---august-2018")
$d.Range(0, 1).Text = ""

# ---------------------------------------------------------------------------
# 3) "Locations considered" bullet - drop "/ M23"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Locations considered: London / Sussex / M23 / remote", $true, $false, $false, $false, $false,
    $true, 1, $false, "Locations considered: London / Sussex / remote", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Perl summary bullet rewrite
# ---------------------------------------------------------------------------
$quote = [char]34
$oldPerl = "Perl (10+ years): " + $quote + "Modern Perl" + $quote + " (Moose/Moo, Plack, Test::*), " + $quote + "the usual" + $quote + " - database apps, web services, data-processing (XML, JSON etc)."
$newPerl = "Perl (10+ years): " + $quote + "Modern Perl" + $quote + " OO (Moose/Moo, Plack) - database apps, web services, data-processing (XML, JSON etc)."
$d.Content.Find.Execute($oldPerl, $true, $false, $false, $false, $false, $true, 1, $false, $newPerl, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Frontend bullet - add "(limited)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Frontend: HTML, CSS, Javascript (JQuery), Bootstrap.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Frontend (limited): HTML, CSS, Javascript (JQuery), Bootstrap.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "Also commercial experience" bullet - capitalise Python/OpenSSL
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Also commercial experience with: python, C/C++, GNU toolset (bash etc), openSSL, XML toolchain (XSLT, XQuery, etc), virtualization.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Also commercial experience with: Python, C/C++, GNU toolset (bash etc), OpenSSL, XML toolchain (XSLT, XQuery, etc), virtualization.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Mendeley role title: "Senior Software Engineer" -> "Java Software Engineer"
#    (first occurrence, paragraph 41 - the Mendeley/Elsevier role heading)
# ---------------------------------------------------------------------------
$pMendeley = $d.Paragraphs.Item(41)
$pMendeley.Range.Find.Execute(
    "Senior Software Engineer", $true, $false, $false, $false, $false,
    $true, 1, $false, "Java Software Engineer", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Mendeley dates: "Jan 2017-present" -> "Jan 2017-July 2018"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "on-site at Mendeley/Elsevier (London / Jan 2017-present)", $true, $false, $false, $false, $false,
    $true, 1, $false, "on-site at Mendeley/Elsevier (London / Jan 2017-July 2018)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9) Mendeley bullet - "product. Mix of" -> "product - mix of"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Reference Manager 2' product. Mix of client-facing and message-processing.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Reference Manager 2' product - mix of client-facing and message-processing.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 10) Mendeley technologies bullet - append ", Payments (Adyen)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Technologies used/learned: Java 8, Dropwizard, Kibana, Redis, TDD, BDD, Oauth2, OpenID Connect, RxJava, Docker/ECS, Terraform, AWS.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Technologies used/learned: Java 8, Dropwizard, Kibana, Redis, TDD, BDD, Oauth2, OpenID Connect, RxJava, Docker/ECS, Terraform, AWS, Payments (Adyen).",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 11) Zoopla role title: "Senior Software Engineer (SEO)" -> "Perl Software Engineer (SEO)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Senior Software Engineer (SEO)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Perl Software Engineer (SEO)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 12) Sky role title: "Senior Software Engineer" -> "Java Software Engineer"
#     (second remaining occurrence, paragraph 49 - the Sky role heading)
# ---------------------------------------------------------------------------
$pSky = $d.Paragraphs.Item(49)
$pSky.Range.Find.Execute(
    "Senior Software Engineer", $true, $false, $false, $false, $false,
    $true, 1, $false, "Java Software Engineer", 2) | Out-Null
